$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '30.471.39'
Set-TextValue $ws.Range('E2') '  +0.43%  '

Set-TextValue $ws.Range('D3') '2.108.54'
Set-TextValue $ws.Range('E3') '  +4.85%  '

Set-TextValue $ws.Range('D4') '1.000'
Set-TextValue $ws.Range('E4') '  -0.20%  '

Set-TextValue $ws.Range('D5') '329.41'
Set-TextValue $ws.Range('E5') '  +1.43%  '

Set-TextValue $ws.Range('D6') '0.9997'
Set-TextValue $ws.Range('E6') '  -0.16%  '

Set-TextValue $ws.Range('D7') '0.5274'
Set-TextValue $ws.Range('E7') '  +2.65%  '

Set-TextValue $ws.Range('D8') '0.4360'
Set-TextValue $ws.Range('E8') '  +1.89%  '

Set-TextValue $ws.Range('D9') '0.08853'
Set-TextValue $ws.Range('E9') '  +1.70%  '

Set-TextValue $ws.Range('D10') '47.10'
Set-TextValue $ws.Range('E10') '  +8.99%  '

Set-TextValue $ws.Range('D11') '1.166'
Set-TextValue $ws.Range('E11') '  +3.00%  '

Set-TextValue $ws.Range('D12') '24.66'
Set-TextValue $ws.Range('E12') '  -0.71%  '

Set-TextValue $ws.Range('D13') '2.099.35'
Set-TextValue $ws.Range('E13') '  +4.57%  '

Set-TextValue $ws.Range('D14') '6.737'
Set-TextValue $ws.Range('E14') '  +2.63%  '

Set-TextValue $ws.Range('D15') '7.796'
Set-TextValue $ws.Range('E15') '  +4.38%  '

Set-TextValue $ws.Range('D16') '96.48'
Set-TextValue $ws.Range('E16') '  +2.22%  '

Set-TextValue $ws.Range('D17') '1.001'
Set-TextValue $ws.Range('E17') '  -0.17%  '

Set-TextValue $ws.Range('D18') '0.00001129'
Set-TextValue $ws.Range('E18') '  +1.48%  '

Set-TextValue $ws.Range('D19') '0.06642'
Set-TextValue $ws.Range('E19') '  +1.83%  '

Set-TextValue $ws.Range('E20') '  +0.64%  '

Set-TextValue $ws.Range('D21') '0.9990'
Set-TextValue $ws.Range('E21') '  -0.23%  '

Set-TextValue $ws.Range('D22') '6.352'
Set-TextValue $ws.Range('E22') '  +2.42%  '

Set-TextValue $ws.Range('D23') '30.523.32'
Set-TextValue $ws.Range('E23') '  +0.40%  '

Set-TextValue $ws.Range('D24') '12.39'
Set-TextValue $ws.Range('E24') '  +4.88%  '

Set-TextValue $ws.Range('D25') '2.332'
Set-TextValue $ws.Range('E25') '  +4.12%  '

Set-TextValue $ws.Range('D26') '2.344.48'
Set-TextValue $ws.Range('E26') '  +4.47%  '

Set-TextValue $ws.Range('D27') '22.44'
Set-TextValue $ws.Range('E27') '  +0.30%  '

Set-TextValue $ws.Range('D28') '2.591'
Set-TextValue $ws.Range('E28') '  +7.35%  '

Set-TextValue $ws.Range('D29') '161.85'
Set-TextValue $ws.Range('E29') '  -0.44%  '

Set-TextValue $ws.Range('D30') '132.55'
Set-TextValue $ws.Range('E30') '  +1.24%  '

Set-TextValue $ws.Range('D31') '1.214'
Set-TextValue $ws.Range('E31') '  +7.00%  '

Set-TextValue $ws.Range('D32') '0.1077'
Set-TextValue $ws.Range('E32') '  +2.24%  '

Set-TextValue $ws.Range('D33') '1.687'
Set-TextValue $ws.Range('E33') '  +23.78%  '

Set-TextValue $ws.Range('D34') '6.196'
Set-TextValue $ws.Range('E34') '  +2.24%  '

Set-TextValue $ws.Range('D35') '3.919'
Set-TextValue $ws.Range('E35') '  +2.13%  '

Set-TextValue $ws.Range('D36') '9.969'
Set-TextValue $ws.Range('E36') '  +9.50%  '

Set-TextValue $ws.Range('D37') '0.02586'
Set-TextValue $ws.Range('E37') '  +2.45%  '

Set-TextValue $ws.Range('D38') '5.505'
Set-TextValue $ws.Range('E38') '  +0.80%  '

Set-TextValue $ws.Range('D39') '0.06713'
Set-TextValue $ws.Range('E39') '  +1.06%  '

Set-TextValue $ws.Range('D40') '12.71'
Set-TextValue $ws.Range('E40') '  +2.62%  '

Set-TextValue $ws.Range('D41') '0.2270'
Set-TextValue $ws.Range('E41') '  +3.78%  '

Set-TextValue $ws.Range('D42') '0.6829'
Set-TextValue $ws.Range('E42') '  +2.82%  '

Set-TextValue $ws.Range('D43') '1.254'
Set-TextValue $ws.Range('E43') '  +1.82%  '

Set-TextValue $ws.Range('B44') 'Frax'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D44') '0.9986'
Set-TextValue $ws.Range('E44') '  -0.25%  '

Set-TextValue $ws.Range('B45') 'EnergySwap'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D45') '14.08'
Set-TextValue $ws.Range('E45') '  +3.53%  '

Set-TextValue $ws.Range('D46') '0.6393'
Set-TextValue $ws.Range('E46') '  +3.75%  '

Set-TextValue $ws.Range('E47') '  +1.54%  '

Set-TextValue $ws.Range('D48') '3.622'
Set-TextValue $ws.Range('E48') '  -1.04%  '

Set-TextValue $ws.Range('D49') '1.251'
Set-TextValue $ws.Range('E49') '  -0.82%  '

Set-TextValue $ws.Range('E50') '  +8.40%  '

Set-TextValue $ws.Range('D51') '82.59'
Set-TextValue $ws.Range('E51') '  +2.56%  '
